$wb = $excel.ActiveWorkbook

# The generated localization report now reflects the files as "Ready for handoff"
# (previously "In Translation"), with refreshed handoff timestamps, and the
# wider "Status"/status-summary columns needed for the longer text.
$newStatus = "Ready for handoff"
$newWidth = 16.3

# Sheet "Overview": per-language status columns (zh-cn / de-de) and the
# "Latest HO Xliff Generate Date" are refreshed.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-15 18:35:42"
$wsOverview.Range("E2:F2").ColumnWidth = $newWidth

# Sheet "zh-cn": status + handoff datetime refreshed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-15 18:35:38"
$wsZhCn.Range("C2").ColumnWidth = $newWidth

# Sheet "de-de": status + handoff datetime refreshed
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-15 18:35:42"
$wsDeDe.Range("C2").ColumnWidth = $newWidth
